# Applies the scheduled-runner update to the Leve profit-calculation sheets.
# For each affected row, refreshes the market-price/profit columns (H:N).
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 107
$ws.Range("H107").Value = 2075.7273
$ws.Range("I107").Value = 1353
$ws.Range("J107").Value = 5328
$ws.Range("K107").Value = 1353
$ws.Range("L107").Value = 5328
$ws.Range("M107").Value = 567
$ws.Range("N107").Value = -9168

# Row 137
$ws.Range("H137").Value = 15718.408
$ws.Range("I137").Value = 18721.053
$ws.Range("J137").Value = 13816.733
$ws.Range("K137").Value = 56163.159
$ws.Range("L137").Value = 41450.199
$ws.Range("M137").Value = -53613.159
$ws.Range("N137").Value = -46550.199

$ws = $wb.Worksheets.Item("ARM")
# Row 32
$ws.Range("H32").Value = 14420.926
$ws.Range("I32").Value = 15359.413
$ws.Range("J32").Value = 2689.8333
$ws.Range("K32").Value = 15359.413
$ws.Range("L32").Value = 2689.8333
$ws.Range("M32").Value = -15072.413
$ws.Range("N32").Value = -3263.8333

# Row 61
$ws.Range("H61").Value = 9465.053
$ws.Range("I61").Value = 1866.5714
$ws.Range("J61").Value = 30740.8
$ws.Range("K61").Value = 1866.5714
$ws.Range("L61").Value = 30740.8
$ws.Range("M61").Value = -1654.5714
$ws.Range("N61").Value = -31164.8

# Row 74
$ws.Range("H74").Value = 150227.05
$ws.Range("I74").Value = 167766.14
$ws.Range("J74").Value = 23945.6
$ws.Range("K74").Value = 167766.14
$ws.Range("L74").Value = 23945.6
$ws.Range("M74").Value = -166892.14
$ws.Range("N74").Value = -25693.6

# Row 77
$ws.Range("H77").Value = 150227.05
$ws.Range("I77").Value = 167766.14
$ws.Range("J77").Value = 23945.6
$ws.Range("K77").Value = 838830.7000000001
$ws.Range("L77").Value = 119728
$ws.Range("M77").Value = -834462.7000000001
$ws.Range("N77").Value = -128464

# Row 110
$ws.Range("H110").Value = 49636.668
$ws.Range("I110").Value = 54847.812
$ws.Range("J110").Value = 7947.5
$ws.Range("K110").Value = 54847.812
$ws.Range("L110").Value = 7947.5
$ws.Range("M110").Value = -52802.812
$ws.Range("N110").Value = -12037.5

# Row 122
$ws.Range("H122").Value = 1406.5
$ws.Range("I122").Value = 1295.32
$ws.Range("J122").Value = 2333
$ws.Range("K122").Value = 3885.96
$ws.Range("L122").Value = 6999
$ws.Range("M122").Value = -1435.96
$ws.Range("N122").Value = -11899

# Row 132
$ws.Range("H132").Value = 2675.1
$ws.Range("I132").Value = 2018.9546
$ws.Range("J132").Value = 4479.5
$ws.Range("K132").Value = 6056.8638
$ws.Range("L132").Value = 13438.5
$ws.Range("M132").Value = -3526.8638
$ws.Range("N132").Value = -18498.5

# Row 136
$ws.Range("H136").Value = 9465.053
$ws.Range("I136").Value = 1866.5714
$ws.Range("J136").Value = 30740.8
$ws.Range("K136").Value = 5599.7142
$ws.Range("L136").Value = 92222.39999999999
$ws.Range("M136").Value = -3049.7142
$ws.Range("N136").Value = -97322.39999999999

$ws = $wb.Worksheets.Item("BSM")
# Row 105
$ws.Range("H105").Value = 2126.9644
$ws.Range("I105").Value = 1929.76
$ws.Range("J105").Value = 3770.3333
$ws.Range("K105").Value = 1929.76
$ws.Range("L105").Value = 3770.3333
$ws.Range("M105").Value = -182.76
$ws.Range("N105").Value = -7264.3333

# Row 107
$ws.Range("H107").Value = 1900.4
$ws.Range("I107").Value = 1996.0278
$ws.Range("J107").Value = 1039.75
$ws.Range("K107").Value = 1996.0278
$ws.Range("L107").Value = 1039.75
$ws.Range("M107").Value = -76.02780000000007
$ws.Range("N107").Value = -4879.75

# Row 134
$ws.Range("H134").Value = 2866.7693
$ws.Range("I134").Value = 2389.85
$ws.Range("J134").Value = 4456.5
$ws.Range("K134").Value = 7169.549999999999
$ws.Range("L134").Value = 13369.5
$ws.Range("M134").Value = -4634.549999999999
$ws.Range("N134").Value = -18439.5

$ws = $wb.Worksheets.Item("CRP")
# Row 31
$ws.Range("H31").Value = 2225384.5
$ws.Range("I31").Value = 3033061
$ws.Range("J31").Value = 4273.75
$ws.Range("K31").Value = 3033061
$ws.Range("L31").Value = 4273.75
$ws.Range("M31").Value = -3032766
$ws.Range("N31").Value = -4863.75

# Row 34
$ws.Range("H34").Value = 2225384.5
$ws.Range("I34").Value = 3033061
$ws.Range("J34").Value = 4273.75
$ws.Range("K34").Value = 3033061
$ws.Range("L34").Value = 4273.75
$ws.Range("M34").Value = -3032859
$ws.Range("N34").Value = -4677.75

# Row 105
$ws.Range("H105").Value = 2777
$ws.Range("I105").Value = 2103.6667
$ws.Range("J105").Value = 3282
$ws.Range("K105").Value = 2103.6667
$ws.Range("L105").Value = 3282
$ws.Range("M105").Value = -356.6667000000002
$ws.Range("N105").Value = -6776

# Row 107
$ws.Range("H107").Value = 660.4
$ws.Range("I107").Value = 566.4
$ws.Range("J107").Value = 848.4
$ws.Range("K107").Value = 566.4
$ws.Range("L107").Value = 848.4
$ws.Range("M107").Value = 1353.6
$ws.Range("N107").Value = -4688.4

# Row 132
$ws.Range("H132").Value = 143391.42
$ws.Range("I132").Value = 143391.42
$ws.Range("J132").Value = 0
$ws.Range("K132").Value = 430174.26
$ws.Range("L132").Value = 0
$ws.Range("M132").Value = -427644.26
$ws.Range("N132").ClearContents()

# Row 134
$ws.Range("H134").Value = 3012.25
$ws.Range("I134").Value = 3012.25
$ws.Range("J134").Value = 0
$ws.Range("K134").Value = 9036.75
$ws.Range("L134").Value = 0
$ws.Range("M134").Value = -6501.75
$ws.Range("N134").ClearContents()

$ws = $wb.Worksheets.Item("CUL")
# Row 107
$ws.Range("H107").Value = 2764.5
$ws.Range("I107").Value = 4539.8
$ws.Range("J107").Value = 1496.4286
$ws.Range("K107").Value = 13619.4
$ws.Range("L107").Value = 4489.2858
$ws.Range("M107").Value = -11699.4
$ws.Range("N107").Value = -8329.2858

$ws = $wb.Worksheets.Item("GSM")
# Row 102
$ws.Range("H102").Value = 30931.295
$ws.Range("I102").Value = 32552
$ws.Range("J102").Value = 5000
$ws.Range("K102").Value = 32552
$ws.Range("L102").Value = 5000
$ws.Range("M102").Value = -30930
$ws.Range("N102").Value = -8244

# Row 132
$ws.Range("H132").Value = 3136.611
$ws.Range("I132").Value = 3061.3635
$ws.Range("J132").Value = 3254.8572
$ws.Range("K132").Value = 9184.0905
$ws.Range("L132").Value = 9764.571599999999
$ws.Range("M132").Value = -6654.0905
$ws.Range("N132").Value = -14824.5716

$ws = $wb.Worksheets.Item("LTW")
# Row 132
$ws.Range("H132").Value = 2575.4348
$ws.Range("I132").Value = 2165.9443
$ws.Range("J132").Value = 4049.6
$ws.Range("K132").Value = 6497.8329
$ws.Range("L132").Value = 12148.8
$ws.Range("M132").Value = -3967.8329
$ws.Range("N132").Value = -17208.8

$ws = $wb.Worksheets.Item("WVR")
# Row 132
$ws.Range("H132").Value = 41953.723
$ws.Range("I132").Value = 52447.145
$ws.Range("J132").Value = 5226.75
$ws.Range("K132").Value = 157341.435
$ws.Range("L132").Value = 15680.25
$ws.Range("M132").Value = -154811.435
$ws.Range("N132").Value = -20740.25

# Row 136
$ws.Range("H136").Value = 25201.357
$ws.Range("I136").Value = 30446.479
$ws.Range("J136").Value = 1073.8
$ws.Range("K136").Value = 91339.43700000001
$ws.Range("L136").Value = 3221.4
$ws.Range("M136").Value = -88789.43700000001
$ws.Range("N136").Value = -8321.4
